# Auto-generated edit script: refresh market-board profit data snapshot
# across all 8 Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 15229.23
$ws.Range("I18").Value = 11999.556
$ws.Range("J18").Value = 22496
$ws.Range("K18").Value = 11999.556
$ws.Range("L18").Value = 22496
$ws.Range("M18").Value = -11715.556
$ws.Range("N18").Value = -23064
$ws.Range("H53").Value = 1346.875
$ws.Range("I53").Value = 2099.4
$ws.Range("J53").Value = 92.666664
$ws.Range("K53").Value = 2099.4
$ws.Range("L53").Value = 92.666664
$ws.Range("M53").Value = -1462.4
$ws.Range("N53").Value = -1366.666664
$ws.Range("H76").Value = 2944
$ws.Range("I76").Value = 2860
$ws.Range("K76").Value = 2860
$ws.Range("M76").Value = -2545
$ws.Range("H79").Value = 2944
$ws.Range("I79").Value = 2860
$ws.Range("K79").Value = 2860
$ws.Range("M79").Value = -1768
$ws.Range("H101").Value = 8888888
$ws.Range("I101").Value = 8888888
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 26666664
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -26665042
$ws.Range("N101").ClearContents()
$ws.Range("H137").Value = 101821.3
$ws.Range("I137").Value = 1801
$ws.Range("K137").Value = 5403
$ws.Range("M137").Value = -2853
$ws.Range("H141").Value = 2547143.5
$ws.Range("I141").Value = 3112002.5
$ws.Range("K141").Value = 9336007.5
$ws.Range("M141").Value = -9330827.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2778777.5
$ws.Range("I2").Value = 5555555
$ws.Range("K2").Value = 5555555
$ws.Range("M2").Value = -5555442
$ws.Range("H32").Value = 5831.2695
$ws.Range("I32").Value = 4102.4556
$ws.Range("K32").Value = 4102.4556
$ws.Range("M32").Value = -3815.4556
$ws.Range("H45").Value = 1336.75
$ws.Range("I45").Value = 1032.8
$ws.Range("K45").Value = 1032.8
$ws.Range("M45").Value = -655.8
$ws.Range("H74").Value = 1130.5
$ws.Range("I74").Value = 607.03845
$ws.Range("K74").Value = 607.03845
$ws.Range("M74").Value = 266.96155
$ws.Range("H77").Value = 1130.5
$ws.Range("I77").Value = 607.03845
$ws.Range("K77").Value = 3035.19225
$ws.Range("M77").Value = 1332.80775
$ws.Range("H102").Value = 1515
$ws.Range("I102").Value = 1503.3334
$ws.Range("J102").Value = 1550
$ws.Range("K102").Value = 1503.3334
$ws.Range("L102").Value = 1550
$ws.Range("M102").Value = 118.6666
$ws.Range("N102").Value = -4794
$ws.Range("H116").Value = 2778777.5
$ws.Range("I116").Value = 5555555
$ws.Range("K116").Value = 5555555
$ws.Range("M116").Value = -5553261
$ws.Range("H132").Value = 2184.5833
$ws.Range("I132").Value = 1569.7059
$ws.Range("K132").Value = 4709.1177
$ws.Range("M132").Value = -2179.1177

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2778777.5
$ws.Range("I3").Value = 5555555
$ws.Range("K3").Value = 5555555
$ws.Range("M3").Value = -5555441
$ws.Range("H99").Value = 1365.5
$ws.Range("I99").Value = 1143.3334
$ws.Range("K99").Value = 1143.3334
$ws.Range("M99").Value = 354.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4350755
$ws.Range("I58").Value = 5437193.5
$ws.Range("K58").Value = 5437193.5
$ws.Range("M58").Value = -5436990.5
$ws.Range("H69").Value = 116600.5
$ws.Range("I69").Value = 13000
$ws.Range("J69").Value = 220201
$ws.Range("K69").Value = 13000
$ws.Range("L69").Value = 220201
$ws.Range("M69").Value = -12251
$ws.Range("N69").Value = -221699
$ws.Range("H72").Value = 116600.5
$ws.Range("I72").Value = 13000
$ws.Range("J72").Value = 220201
$ws.Range("K72").Value = 39000
$ws.Range("L72").Value = 660603
$ws.Range("M72").Value = -35256
$ws.Range("N72").Value = -668091
$ws.Range("H86").Value = 3019.4
$ws.Range("I86").Value = 2649
$ws.Range("K86").Value = 2649
$ws.Range("M86").Value = -1526
$ws.Range("H89").Value = 3019.4
$ws.Range("I89").Value = 2649
$ws.Range("K89").Value = 13245
$ws.Range("M89").Value = -7629
$ws.Range("H122").Value = 3160.2727
$ws.Range("I122").Value = 1917.4
$ws.Range("J122").Value = 4196
$ws.Range("K122").Value = 5752.200000000001
$ws.Range("L122").Value = 12588
$ws.Range("M122").Value = -3302.200000000001
$ws.Range("N122").Value = -17488
$ws.Range("H132").Value = 1346.6364
$ws.Range("I132").Value = 953.3684
$ws.Range("K132").Value = 2860.1052
$ws.Range("M132").Value = -330.1052
$ws.Range("H134").Value = 2751.353
$ws.Range("I134").Value = 2537
$ws.Range("K134").Value = 7611
$ws.Range("M134").Value = -5076
$ws.Range("H136").Value = 4350755
$ws.Range("I136").Value = 5437193.5
$ws.Range("K136").Value = 16311580.5
$ws.Range("M136").Value = -16309030.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H113").Value = 383432.22
$ws.Range("J113").Value = 400681.2
$ws.Range("L113").Value = 1202043.6
$ws.Range("N113").Value = -1206383.6
$ws.Range("H131").Value = 13433.063
$ws.Range("J131").Value = 13855.787
$ws.Range("L131").Value = 41567.361
$ws.Range("N131").Value = -51647.361
$ws.Range("H136").Value = 1380.3334
$ws.Range("I136").Value = 1380.3334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4141.0002
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 958.9997999999996
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 7135.421
$ws.Range("I139").Value = 7420.722
$ws.Range("K139").Value = 22262.166
$ws.Range("M139").Value = -17122.166
$ws.Range("H140").Value = 2260.2222
$ws.Range("I140").Value = 1419.2
$ws.Range("K140").Value = 4257.6
$ws.Range("M140").Value = 922.3999999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1224.0714
$ws.Range("J113").Value = 1498
$ws.Range("L113").Value = 1498
$ws.Range("N113").Value = -5838
$ws.Range("H132").Value = 3850434.5
$ws.Range("J132").Value = 5696.5
$ws.Range("L132").Value = 17089.5
$ws.Range("N132").Value = -22149.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2704.8667
$ws.Range("I7").Value = 2762.75
$ws.Range("K7").Value = 2762.75
$ws.Range("M7").Value = -2650.75
$ws.Range("H16").Value = 4969.143
$ws.Range("I16").Value = 6188.273
$ws.Range("J16").Value = 499
$ws.Range("K16").Value = 6188.273
$ws.Range("L16").Value = 499
$ws.Range("M16").Value = -6018.273
$ws.Range("N16").Value = -839
$ws.Range("H22").Value = 1252.6
$ws.Range("I22").Value = 737.875
$ws.Range("J22").Value = 1595.75
$ws.Range("K22").Value = 737.875
$ws.Range("L22").Value = 1595.75
$ws.Range("M22").Value = -442.875
$ws.Range("N22").Value = -2185.75
$ws.Range("H27").Value = 1252.6
$ws.Range("I27").Value = 737.875
$ws.Range("J27").Value = 1595.75
$ws.Range("K27").Value = 737.875
$ws.Range("L27").Value = 1595.75
$ws.Range("M27").Value = -630.875
$ws.Range("N27").Value = -1809.75
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 425.9091
$ws.Range("I55").Value = 370.36365
$ws.Range("J55").Value = 481.45456
$ws.Range("K55").Value = 370.36365
$ws.Range("L55").Value = 481.45456
$ws.Range("M55").Value = -197.36365
$ws.Range("N55").Value = -827.45456
$ws.Range("H61").Value = 4577.4287
$ws.Range("J61").Value = 4755.5713
$ws.Range("L61").Value = 4755.5713
$ws.Range("N61").Value = -5159.5713
$ws.Range("H113").Value = 4577.4287
$ws.Range("J113").Value = 4755.5713
$ws.Range("L113").Value = 4755.5713
$ws.Range("N113").Value = -9095.5713
$ws.Range("H126").Value = 2704.8667
$ws.Range("I126").Value = 2762.75
$ws.Range("K126").Value = 8288.25
$ws.Range("M126").Value = -5818.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 40518.05
$ws.Range("I122").Value = 50258.375
$ws.Range("J122").Value = 1556.75
$ws.Range("K122").Value = 150775.125
$ws.Range("L122").Value = 4670.25
$ws.Range("M122").Value = -148325.125
$ws.Range("N122").Value = -9570.25
$ws.Range("H132").Value = 847.28815
$ws.Range("I132").Value = 720
$ws.Range("J132").Value = 1345.8334
$ws.Range("K132").Value = 720
$ws.Range("L132").Value = 4037.5002
$ws.Range("M132").Value = 370
$ws.Range("N132").Value = -9097.5002
